# The workbook contains a weekly-updated price history table (rows 2-64,
# header in row 1). This edit adds one new week's record at the top of the
# data block (row 12) and appends the newest week's record at the bottom
# (new row 66), pushing every data row from 12..64 down by one.
#
# Row 12 is where the new record is inserted; everything that used to be
# in rows 12..64 shifts down to 13..65 automatically via Rows.Insert().
# Then a brand-new row 66 is appended after the (shifted) former last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at position 12 - shifts old rows 12..64 -> 13..65
$ws.Rows.Item(12).Insert()

# 2) Populate the newly inserted row 12 with this week's record.
#    Columns that never vary across this table (A,B,C,E,F,G,H,I,N,O,Q,R)
#    are copied from the constant pattern used by every other row.
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 44831
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 100114007
$ws.Cells.Item(12, 7).Value = "Jengibre"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 600
$ws.Cells.Item(12, 11).Value = 14000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 14500
$ws.Cells.Item(12, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(12, 15).Value = 'Perú'
$ws.Cells.Item(12, 16).Value = 1115
$ws.Cells.Item(12, 17).Value = 13
$ws.Cells.Item(12, 18).Value = "Hortaliza"

# 3) Append a brand-new row 66 (one row after the shifted former last row,
#    which is now row 65) with the latest week's record.
$ws.Cells.Item(66, 1).Value = 8
$ws.Cells.Item(66, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(66, 3).Value = "Coquimbo"
$ws.Cells.Item(66, 4).Value = 44832
$ws.Cells.Item(66, 5).Value = 4
$ws.Cells.Item(66, 6).Value = 100114007
$ws.Cells.Item(66, 7).Value = "Jengibre"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 540
$ws.Cells.Item(66, 11).Value = 14000
$ws.Cells.Item(66, 12).Value = 15000
$ws.Cells.Item(66, 13).Value = 14500
$ws.Cells.Item(66, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(66, 15).Value = 'Perú'
$ws.Cells.Item(66, 16).Value = 1115
$ws.Cells.Item(66, 17).Value = 13
$ws.Cells.Item(66, 18).Value = "Hortaliza"

# 4) Give the date cells (column D) the same date-number-format style (s="2")
#    used by every other row in the table, same as D12/D66 siblings.
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(13, 4).NumberFormat
$ws.Cells.Item(66, 4).NumberFormat = $ws.Cells.Item(65, 4).NumberFormat
